# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
# (crypto price/volume refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several new Price values parse as plain numbers (e.g. "19.45"), but the
# sheet stores this column as plain text (the diff keeps t="inlineStr").
# Force each such cell to Text format *before* writing so Excel's COM type
# inference doesn't silently convert it to a numeric cell.
foreach ($cellRef in @(
        'D5',
        'D8',
        'D10',
        'D15',
        'D16',
        'D19',
        'D21',
        'D22',
        'D24',
        'D25',
        'D27',
        'D28',
        'D29',
        'D30',
        'D32',
        'D33',
        'D39',
        'D42',
        'D43',
        'D44',
        'D46',
        'D47',
        'D49',
        'D50',
        'D51'
    )) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.624.62'
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').Value = '1.596.34'
$ws.Range('E3').Value = '  +0.42%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '211.53'
$ws.Range('E5').Value = '  +0.25%  '

$ws.Range('E6').Value = '  +1.21%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '0.0617'
$ws.Range('E8').Value = '  +0.34%  '

$ws.Range('E9').Value = '  -0.38%  '

$ws.Range('D10').Value = '19.45'
$ws.Range('E10').Value = '  -0.67%  '

$ws.Range('E11').Value = '  +0.44%  '

$ws.Range('D12').Value = '1.819.69'
$ws.Range('E12').Value = '  +0.46%  '

$ws.Range('D13').Value = '1.590.70'
$ws.Range('E13').Value = '  +0.09%  '

$ws.Range('E14').Value = '  +0.24%  '

$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  -0.14%  '

$ws.Range('D16').Value = '64.61'
$ws.Range('E16').Value = '  -0.37%  '

$ws.Range('D17').Value = '26.606.90'
$ws.Range('E17').Value = '  -0.01%  '

$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  +0.45%  '

$ws.Range('D19').Value = '208.72'
$ws.Range('E19').Value = '  +0.30%  '

$ws.Range('E20').Value = '  -0.02%  '

$ws.Range('D21').Value = '6.97'
$ws.Range('E21').Value = '  +3.60%  '

$ws.Range('D22').Value = '4.26'
$ws.Range('E22').Value = '  +0.43%  '

$ws.Range('E23').Value = '  -1.66%  '

$ws.Range('D24').Value = '8.88'
$ws.Range('E24').Value = '  -0.03%  '

$ws.Range('D25').Value = '145.27'
$ws.Range('E25').Value = '  -1.06%  '

$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('D27').Value = '7.10'
$ws.Range('E27').Value = '  -1.96%  '

$ws.Range('D28').Value = '0.115'
$ws.Range('E28').Value = '  +0.86%  '

$ws.Range('D29').Value = '15.24'
$ws.Range('E29').Value = '  -0.24%  '

$ws.Range('D30').Value = '0.0506'
$ws.Range('E30').Value = '  -0.02%  '

$ws.Range('E31').Value = '  +0.50%  '

$ws.Range('D32').Value = '3.23'
$ws.Range('E32').Value = '  -0.08%  '

$ws.Range('D33').Value = '0.653'
$ws.Range('E33').Value = '  -1.60%  '

$ws.Range('E34').Value = '  +0.96%  '

$ws.Range('D35').Value = '1.283.56'
$ws.Range('E35').Value = '  -1.82%  '

$ws.Range('E36').Value = '  +1.09%  '

$ws.Range('E37').Value = '  +0.38%  '

$ws.Range('E38').Value = '  -0.44%  '

$ws.Range('D39').Value = '0.843'
$ws.Range('E39').Value = '  +1.72%  '

$ws.Range('E40').Value = '  +0.06%  '

$ws.Range('E41').Value = '  +2.08%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.785'
$ws.Range('E42').Value = '  -0.89%  '

$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.19'
$ws.Range('E43').Value = '  +1.20%  '

$ws.Range('D44').Value = '64.26'
$ws.Range('E44').Value = '  +2.27%  '

$ws.Range('D45').Value = '1.732.29'
$ws.Range('E45').Value = '  +0.47%  '

$ws.Range('D46').Value = '0.911'
$ws.Range('E46').Value = '  +8.83%  '

$ws.Range('D47').Value = '89.63'

$ws.Range('E48').Value = '  -0.84%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.102'
$ws.Range('E49').Value = '  +4.31%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0506'
$ws.Range('E50').Value = '  +0.35%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.49'
$ws.Range('E51').Value = '  -0.75%  '

